$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for all data rows (2-506)
# from serial date 45192 to 45202.
$ws.Range("C2:C506").Value = 45202
